$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-257 all change from 45190 to 45192 (date serial)
$ws.Range("C2:C257").Value = 45192
